$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 2331
$ws.Range("J3").Value = 2416
$ws.Range("E4").Value = 1986
$ws.Range("F4").Value = 1881
$ws.Range("J4").Value = 548
$ws.Range("J5").Value = 171
$ws.Range("I6").Value = 8968
$ws.Range("J6").Value = 3027
$ws.Range("E7").Value = 25990
$ws.Range("F7").Value = 24071
$ws.Range("I7").Value = 26202
$ws.Range("J7").Value = 8493

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 66
$ws.Range("J4").Value = 36
$ws.Range("J7").Value = 263
$ws.Range("J8").Value = 537
$ws.Range("J9").Value = 54
$ws.Range("J10").Value = 52
$ws.Range("J11").Value = 117
$ws.Range("J14").Value = 31
$ws.Range("J15").Value = 104
$ws.Range("J22").Value = 17
$ws.Range("J23").Value = 83
$ws.Range("J25").Value = 51
$ws.Range("J29").Value = 476
$ws.Range("J31").Value = 65
$ws.Range("J33").Value = 351
$ws.Range("J36").Value = 131
$ws.Range("I37").Value = 807
$ws.Range("J37").Value = 289
$ws.Range("J43").Value = 79
$ws.Range("J48").Value = 81
$ws.Range("J51").Value = 114
$ws.Range("J52").Value = 208
$ws.Range("J53").Value = 82
$ws.Range("J54").Value = 171
$ws.Range("J55").Value = 98
$ws.Range("J60").Value = 56
$ws.Range("J61").Value = 16
$ws.Range("E63").Value = 331
$ws.Range("F63").Value = 173
$ws.Range("J63").Value = 34
$ws.Range("J65").Value = 222
$ws.Range("J66").Value = 22
$ws.Range("J67").Value = 309
$ws.Range("J77").Value = 63
$ws.Range("J78").Value = 118
$ws.Range("J79").Value = 259
$ws.Range("J83").Value = 203
$ws.Range("J84").Value = 82
$ws.Range("J85").Value = 399
$ws.Range("J86").Value = 47
$ws.Range("J89").Value = 89
$ws.Range("J90").Value = 95
$ws.Range("J91").Value = 94
$ws.Range("J94").Value = 71
$ws.Range("J95").Value = 126
$ws.Range("J96").Value = 97
$ws.Range("J99").Value = 117
$ws.Range("E101").Value = 25990
$ws.Range("F101").Value = 24071
$ws.Range("I101").Value = 26202
$ws.Range("J101").Value = 8493

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J3").Value = 151
$ws.Range("J7").Value = 399

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J2").Value = 50
$ws.Range("J3").Value = 60
$ws.Range("J7").Value = 208

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J3").Value = 22
$ws.Range("J7").Value = 117

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J3").Value = 174
$ws.Range("J6").Value = 157
$ws.Range("J7").Value = 537

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J2").Value = 15
$ws.Range("J6").Value = 47
$ws.Range("J7").Value = 82

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J3").Value = 78
$ws.Range("J6").Value = 91
$ws.Range("J7").Value = 263

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J6").Value = 30
$ws.Range("J7").Value = 89

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J3").Value = 29
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 97

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 31

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 86
$ws.Range("I6").Value = 244
$ws.Range("I7").Value = 807
$ws.Range("J7").Value = 289

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 39
$ws.Range("J7").Value = 117

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 125
$ws.Range("J7").Value = 309

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J2").Value = 25
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 82

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J6").Value = 81
$ws.Range("J7").Value = 222

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 62
$ws.Range("J3").Value = 71
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 203

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 46
$ws.Range("J7").Value = 126

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J4").Value = 19
$ws.Range("J6").Value = 124
$ws.Range("J7").Value = 351

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J3").Value = 34
$ws.Range("J7").Value = 171

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 140
$ws.Range("J5").Value = 18
$ws.Range("J6").Value = 130
$ws.Range("J7").Value = 476

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J2").Value = 17
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J2").Value = 30
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 118

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J2").Value = 26
$ws.Range("J7").Value = 98

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J2").Value = 23
$ws.Range("J3").Value = 29
$ws.Range("J4").Value = 7
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 83

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 94

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J3").Value = 98
$ws.Range("J7").Value = 259

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 131

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J2").Value = 22
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J2").Value = 29
$ws.Range("J7").Value = 104

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 22

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J2").Value = 14
$ws.Range("J7").Value = 54

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 66

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 30
$ws.Range("J7").Value = 95

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 30
$ws.Range("J3").Value = 33
$ws.Range("J7").Value = 114

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J6").Value = 19
$ws.Range("J7").Value = 56

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J6").Value = 49
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("J2").Value = 9
$ws.Range("J7").Value = 17

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 63

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 36

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("J2").Value = 6
$ws.Range("J7").Value = 16
